$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.423324942588806
$ws.Range("B1").Value = 1.650299787521362
$ws.Range("C1").Value = 2.183262825012207
$ws.Range("D1").Value = 2.040988922119141
$ws.Range("E1").Value = 1.327208757400513
